# Regenerate orders with updated distance/size labels.
# D64 -> D69, D51 -> D55, D80 -> D86, S30 -> S31 (S20/S25 unchanged)
# across every cell of the sheet (Condition, Filename_Left, Filename_Right,
# Distance and Size columns all embed these substrings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D64", "D69") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("D80", "D86") | Out-Null
$used.Replace("S30", "S31") | Out-Null
